# Auto-generated edit script: update crypto price/volume table
# Applies the GitHub Actions "Updated cryptos list" data refresh diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "27.817.55"
$ws.Range("E2").Value = "  +1.00%  "

# Row 3
$ws.Range("D3").Value = "1.756.17"
$ws.Range("E3").Value = "  +0.12%  "

# Row 4
$ws.Range("E4").Value = "  +0.07%  "

# Row 5
$ws.Range("D5").Value = "'328.05"
$ws.Range("E5").Value = "  +1.15%  "

# Row 7
$ws.Range("D7").Value = "'0.4575"
$ws.Range("E7").Value = "  +0.30%  "

# Row 8
$ws.Range("D8").Value = "'0.3491"
$ws.Range("E8").Value = "  -1.77%  "

# Row 9
$ws.Range("D9").Value = "'41.91"
$ws.Range("E9").Value = "  +0.95%  "

# Row 10
$ws.Range("D10").Value = "'0.07346"
$ws.Range("E10").Value = "  -1.63%  "

# Row 11
$ws.Range("D11").Value = "'1.083"
$ws.Range("E11").Value = "  -0.35%  "

# Row 12
$ws.Range("D12").Value = "'1.002"
$ws.Range("E12").Value = "  +0.10%  "

# Row 13
$ws.Range("D13").Value = "'20.60"
$ws.Range("E13").Value = "  -0.78%  "

# Row 14
$ws.Range("D14").Value = "'5.975"
$ws.Range("E14").Value = "  -0.68%  "

# Row 15
$ws.Range("D15").Value = "'7.169"
$ws.Range("E15").Value = "  -0.09%  "

# Row 16
$ws.Range("D16").Value = "1.756.46"
$ws.Range("E16").Value = "  -0.10%  "

# Row 17
$ws.Range("D17").Value = "'91.61"
$ws.Range("E17").Value = "  -2.60%  "

# Row 18
$ws.Range("E18").Value = "  -0.31%  "

# Row 19
$ws.Range("D19").Value = "'0.06411"
$ws.Range("E19").Value = "  +0.10%  "

# Row 20
$ws.Range("E20").Value = "  +0.07%  "

# Row 21
$ws.Range("D21").Value = "'16.83"
$ws.Range("E21").Value = "  -1.72%  "

# Row 22
$ws.Range("D22").Value = "'5.736"
$ws.Range("E22").Value = "  -0.21%  "

# Row 23
$ws.Range("D23").Value = "27.859.21"
$ws.Range("E23").Value = "  +0.94%  "

# Row 24
$ws.Range("E24").Value = "  -0.41%  "

# Row 25
$ws.Range("D25").Value = "'2.159"
$ws.Range("E25").Value = "  +3.68%  "

# Row 26
$ws.Range("D26").Value = "'162.49"
$ws.Range("E26").Value = "  -1.99%  "

# Row 27
$ws.Range("D27").Value = "'20.01"
$ws.Range("E27").Value = "  -0.59%  "

# Row 28
$ws.Range("D28").Value = "1.959.61"
$ws.Range("E28").Value = "  +0.08%  "

# Row 29
$ws.Range("D29").Value = "'2.154"
$ws.Range("E29").Value = "  +1.29%  "

# Row 30
$ws.Range("D30").Value = "'123.12"
$ws.Range("E30").Value = "  -2.07%  "

# Row 31
$ws.Range("D31").Value = "'1.081"
$ws.Range("E31").Value = "  -0.53%  "

# Row 32
$ws.Range("D32").Value = "'0.09276"
$ws.Range("E32").Value = "  +0.52%  "

# Row 33
$ws.Range("D33").Value = "'3.641"
$ws.Range("E33").Value = "  -0.43%  "

# Row 34
$ws.Range("D34").Value = "'5.520"
$ws.Range("E34").Value = "  -0.16%  "

# Row 35
$ws.Range("D35").Value = "'11.72"
$ws.Range("E35").Value = "  -0.05%  "

# Row 36
$ws.Range("D36").Value = "'0.06100"
$ws.Range("E36").Value = "  +1.32%  "

# Row 37
$ws.Range("E37").Value = "  -1.51%  "

# Row 38
$ws.Range("D38").Value = "'0.2062"
$ws.Range("E38").Value = "  -1.35%  "

# Row 39
$ws.Range("D39").Value = "'4.886"
$ws.Range("E39").Value = "  -0.77%  "

# Row 40
$ws.Range("D40").Value = "'0.6185"
$ws.Range("E40").Value = "  -1.79%  "

# Row 41
$ws.Range("D41").Value = "'1.178"
$ws.Range("E41").Value = "  -0.17%  "

# Row 42
$ws.Range("E42").Value = "  -1.50%  "

# Row 43
$ws.Range("D43").Value = "'7.775"
$ws.Range("E43").Value = "  -0.34%  "

# Row 44
$ws.Range("D44").Value = "'13.13"
$ws.Range("E44").Value = "  -0.72%  "

# Row 45
$ws.Range("D45").Value = "'3.719"
$ws.Range("E45").Value = "  +0.07%  "

# Row 46
$ws.Range("D46").Value = "'0.5795"
$ws.Range("E46").Value = "  -1.04%  "

# Row 47
$ws.Range("D47").Value = "'122.21"
$ws.Range("E47").Value = "  +0.14%  "

# Row 48
$ws.Range("D48").Value = "'1.923"
$ws.Range("E48").Value = "  -0.56%  "

# Row 49
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "'0.06781"
$ws.Range("E49").Value = "  -1.62%  "

# Row 50
$ws.Range("B50").Value = "EOS"
$ws.Range("C50").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D50").Value = "'1.119"
$ws.Range("E50").Value = "  -0.96%  "

# Row 51
$ws.Range("D51").Value = "'72.13"
$ws.Range("E51").Value = "  +0.10%  "
